{"js": "// Apply Swahili (Kenya) translations to the subtitle document.\n// Each entry is an exact [originalText, translatedText] pair taken from\n// the document; every originalText occurs exactly once in the body, so a\n// simple search + replace is unambiguous for each run.\nconst pairs = [\n  [\n    \"Format has been corrected not the timing\",\n    \"Umbizo limesahihishwa sio wakati\"\n  ],\n  [\n    \"I added 25 seconds to each timing to correct for the intro song -john argentino\",\n    \"Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino\"\n  ],\n  [\n    \"The airport problem - subtitles:\",\n    \"Tatizo la uwanja wa ndege - manukuu:\"\n  ],\n  [\n    \"The administrations of three\",\n    \"Utawala wa tatu\"\n  ],\n  [\n    \"neighboring cities: A, B and C decided\",\n    \"miji jirani: A, B na C waliamua\"\n  ],\n  [\n    \"to build an airport dividing the costs of\",\n    \"kujenga uwanja wa ndege unaogawanya gharama za\"\n  ],\n  [\n    \"implementation. The condition on the\",\n    \"utekelezaji. Hali juu ya\"\n  ],\n  [\n    \"choice of the most suitable place is\",\n    \"uchaguzi wa mahali pa kufaa zaidi ni\"\n  ],\n  [\n    \"that the sum of the distances from each\",\n    \"kwamba jumla ya umbali kutoka kwa kila mmoja\"\n  ],\n  [\n    \"city to the airport is as small as\",\n    \"mji kwa uwanja wa ndege ni ndogo kama\"\n  ],\n  [\n    \"possible. The team of experts in charge\",\n    \"inawezekana. Timu ya wataalam wanaohusika\"\n  ],\n  [\n    \"of the work has created a model to get\",\n    \"ya kazi imeunda mfano wa kupata\"\n  ],\n  [\n    \"a preliminary idea of where to place the\",\n    \"wazo la awali la mahali pa kuweka\"\n  ],\n  [\n    \"structure. At their disposal there are\",\n    \"muundo. Ovyo wao wapo\"\n  ],\n  [\n    \"some snails a big metal ring and a long\",\n    \"konokono wengine pete kubwa ya chuma na ndefu\"\n  ],\n  [\n    \"string.\",\n    \"kamba.\"\n  ],\n  [\n    \"Explain how the team can manage to use\",\n    \"Eleza jinsi timu inaweza kusimamia matumizi\"\n  ],\n  [\n    \"the materials to tell approximately the\",\n    \"nyenzo za kusema takriban\"\n  ],\n  [\n    \"ideal location of the airport. Imagine\",\n    \"eneo bora la uwanja wa ndege. Fikiria\"\n  ],\n  [\n    \"that the cities are placed at the\",\n    \"kwamba miji imewekwa kwenye\"\n  ],\n  [\n    \"vertices of a triangle which is\",\n    \"vipeo vya pembetatu ambayo ni\"\n  ],\n  [\n    \"obviously reproduced in scale as\",\n    \"kwa hakika imetolewa tena kwa kiwango kama\"\n  ],\n  [\n    \"shown in figure. This is one possible\",\n    \"inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana\"\n  ],\n  [\n    \"setting the rope starts from one nail,\",\n    \"kuweka kamba huanza kutoka msumari mmoja,\"\n  ],\n  [\n    \"goes inside the ring, goes around the\",\n    \"huenda ndani ya pete, huzunguka\"\n  ],\n  [\n    \"other nail, the third nail, inside the\",\n    \"msumari mwingine, msumari wa tatu, ndani ya\"\n  ],\n  [\n    \"ring again and now you can just pull the\",\n    \"pete tena na sasa unaweza kuvuta tu\"\n  ],\n  [\n    \"rope in order to find the point that\",\n    \"kamba ili kupata uhakika huo\"\n  ],\n  [\n    \"you're looking for. In order to reach the\",\n    \"unatafuta. Ili kufikia\"\n  ],\n  [\n    \"point, we have to move the rope a bit\",\n    \"uhakika, tunapaswa kusonga kamba kidogo\"\n  ],\n  [\n    \"because there is some \",\n    \"kwa sababu kuna \"\n  ],\n  [\n    \"resistance\",\n    \"upinzani\"\n  ],\n  [\n    \" caused\",\n    \" uliosababishwa\"\n  ],\n  [\n    \"by the materials that we are using but\",\n    \"kwa nyenzo ambazo tunatumia lakini\"\n  ],\n  [\n    \"after a while you'll reach a position from\",\n    \"baada ya muda utafikia nafasi kutoka\"\n  ],\n  [\n    \"which the ring doesn't move anymore,\",\n    \"ambayo pete haisogei tena,\"\n  ],\n  [\n    \"which is more or less this one. And as\",\n    \"ambayo ni zaidi au chini ya hii. Na kama\"\n  ],\n  [\n    \"between the ring and the nails are\",\n    \"kati ya pete na misumari ni\"\n  ],\n  [\n    \"placed more or less 120 degrees from one\",\n    \"kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja\"\n  ],\n  [\n    \"another which is 1/3 of a circumference,\",\n    \"nyingine ambayo ni 1/3 ya mduara,\"\n  ],\n  [\n    \"and that's the point that we're looking\",\n    \"na hiyo ndiyo hatua tunayoiangalia\"\n  ],\n  [\n    \"for: the minimum distance between the\",\n    \"kwa: umbali wa chini kati ya\"\n  ],\n  [\n    \"nails and the airport when you sum it\",\n    \"misumari na uwanja wa ndege unapojumlisha\"\n  ],\n  [\n    \"ogether\",\n    \"pamoja\"\n  ],\n  [\n    \"[Music]\",\n    \"[Muziki]\"\n  ]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply Swahili (Kenya) translations to the subtitle document.\n# Each entry is an exact (originalText, translatedText) pair taken from\n# the document; every originalText occurs exactly once in the body, so\n# Find/Replace with MatchCase and Replace=wdReplaceOne is unambiguous.\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$pairs = @(\n  ,@('Format has been corrected not the timing', 'Umbizo limesahihishwa sio wakati')\n  ,@('I added 25 seconds to each timing to correct for the intro song -john argentino', 'Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino')\n  ,@('The airport problem - subtitles:', 'Tatizo la uwanja wa ndege - manukuu:')\n  ,@('The administrations of three', 'Utawala wa tatu')\n  ,@('neighboring cities: A, B and C decided', 'miji jirani: A, B na C waliamua')\n  ,@('to build an airport dividing the costs of', 'kujenga uwanja wa ndege unaogawanya gharama za')\n  ,@('implementation. The condition on the', 'utekelezaji. Hali juu ya')\n  ,@('choice of the most suitable place is', 'uchaguzi wa mahali pa kufaa zaidi ni')\n  ,@('that the sum of the distances from each', 'kwamba jumla ya umbali kutoka kwa kila mmoja')\n  ,@('city to the airport is as small as', 'mji kwa uwanja wa ndege ni ndogo kama')\n  ,@('possible. The team of experts in charge', 'inawezekana. Timu ya wataalam wanaohusika')\n  ,@('of the work has created a model to get', 'ya kazi imeunda mfano wa kupata')\n  ,@('a preliminary idea of where to place the', 'wazo la awali la mahali pa kuweka')\n  ,@('structure. At their disposal there are', 'muundo. Ovyo wao wapo')\n  ,@('some snails a big metal ring and a long', 'konokono wengine pete kubwa ya chuma na ndefu')\n  ,@('string.', 'kamba.')\n  ,@('Explain how the team can manage to use', 'Eleza jinsi timu inaweza kusimamia matumizi')\n  ,@('the materials to tell approximately the', 'nyenzo za kusema takriban')\n  ,@('ideal location of the airport. Imagine', 'eneo bora la uwanja wa ndege. Fikiria')\n  ,@('that the cities are placed at the', 'kwamba miji imewekwa kwenye')\n  ,@('vertices of a triangle which is', 'vipeo vya pembetatu ambayo ni')\n  ,@('obviously reproduced in scale as', 'kwa hakika imetolewa tena kwa kiwango kama')\n  ,@('shown in figure. This is one possible', 'inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana')\n  ,@('setting the rope starts from one nail,', 'kuweka kamba huanza kutoka msumari mmoja,')\n  ,@('goes inside the ring, goes around the', 'huenda ndani ya pete, huzunguka')\n  ,@('other nail, the third nail, inside the', 'msumari mwingine, msumari wa tatu, ndani ya')\n  ,@('ring again and now you can just pull the', 'pete tena na sasa unaweza kuvuta tu')\n  ,@('rope in order to find the point that', 'kamba ili kupata uhakika huo')\n  ,@('you''re looking for. In order to reach the', 'unatafuta. Ili kufikia')\n  ,@('point, we have to move the rope a bit', 'uhakika, tunapaswa kusonga kamba kidogo')\n  ,@('because there is some ', 'kwa sababu kuna ')\n  ,@('resistance', 'upinzani')\n  ,@(' caused', ' uliosababishwa')\n  ,@('by the materials that we are using but', 'kwa nyenzo ambazo tunatumia lakini')\n  ,@('after a while you''ll reach a position from', 'baada ya muda utafikia nafasi kutoka')\n  ,@('which the ring doesn''t move anymore,', 'ambayo pete haisogei tena,')\n  ,@('which is more or less this one. And as', 'ambayo ni zaidi au chini ya hii. Na kama')\n  ,@('between the ring and the nails are', 'kati ya pete na misumari ni')\n  ,@('placed more or less 120 degrees from one', 'kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja')\n  ,@('another which is 1/3 of a circumference,', 'nyingine ambayo ni 1/3 ya mduara,')\n  ,@('and that''s the point that we''re looking', 'na hiyo ndiyo hatua tunayoiangalia')\n  ,@('for: the minimum distance between the', 'kwa: umbali wa chini kati ya')\n  ,@('nails and the airport when you sum it', 'misumari na uwanja wa ndege unapojumlisha')\n  ,@('ogether', 'pamoja')\n  ,@('[Music]', '[Muziki]')\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n"}
